$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear previously used range contents (keep formatting) so stale values in shifted columns do not linger
$ws.Range("A1:J60").ClearContents()

# Header row: B1:I1 = d=1,2,3,4,5,6,7,10
$ws.Range("B1").Value = "d=1"
$ws.Range("C1").Value = "d=2"
$ws.Range("D1").Value = "d=3"
$ws.Range("E1").Value = "d=4"
$ws.Range("F1").Value = "d=5"
$ws.Range("G1").Value = "d=6"
$ws.Range("H1").Value = "d=7"
$ws.Range("I1").Value = "d=10"

$ws.Range("A2").Value = "ARMA_I(0,1,0)"
$ws.Range("B2").Value = 82.40641713471288
$ws.Range("A3").Value = "ARMA_I(0,1,1)"
$ws.Range("B3").Value = 86.73939412659749
$ws.Range("A4").Value = "ARMA_I(0,1,2)"
$ws.Range("B4").Value = 89.80045417174072
$ws.Range("A5").Value = "ARMA_I(0,10,0)"
$ws.Range("I5").Value = 95.88925451112898
$ws.Range("A6").Value = "ARMA_I(0,10,1)"
$ws.Range("I6").Value = 95.88701720109172
$ws.Range("A7").Value = "ARMA_I(0,10,2)"
$ws.Range("I7").Value = 95.87308791677911
$ws.Range("A8").Value = "ARMA_I(0,2,0)"
$ws.Range("C8").Value = 98.41432840090805
$ws.Range("A9").Value = "ARMA_I(0,2,1)"
$ws.Range("C9").Value = 97.0255619905421
$ws.Range("A10").Value = "ARMA_I(0,2,2)"
$ws.Range("C10").Value = 98.11787034460517
$ws.Range("A11").Value = "ARMA_I(0,3,0)"
$ws.Range("D11").Value = 99.32639786170625
$ws.Range("A12").Value = "ARMA_I(0,3,1)"
$ws.Range("D12").Value = 99.21500274414939
$ws.Range("A13").Value = "ARMA_I(0,3,2)"
$ws.Range("D13").Value = 99.27432671514555
$ws.Range("A14").Value = "ARMA_I(0,4,0)"
$ws.Range("E14").Value = 98.72743293293674
$ws.Range("A15").Value = "ARMA_I(0,4,1)"
$ws.Range("E15").Value = 98.79197456695178
$ws.Range("A16").Value = "ARMA_I(0,4,2)"
$ws.Range("E16").Value = 98.7566203622452
$ws.Range("A17").Value = "ARMA_I(0,5,0)"
$ws.Range("F17").Value = 98.2289834120438
$ws.Range("A18").Value = "ARMA_I(0,5,1)"
$ws.Range("F18").Value = 98.29785025341207
$ws.Range("A19").Value = "ARMA_I(0,5,2)"
$ws.Range("F19").Value = 98.21321661171832
$ws.Range("A20").Value = "ARMA_I(0,6,0)"
$ws.Range("G20").Value = 97.81703327887682
$ws.Range("A21").Value = "ARMA_I(0,6,1)"
$ws.Range("G21").Value = 97.852526539316
$ws.Range("A22").Value = "ARMA_I(0,6,2)"
$ws.Range("G22").Value = 97.7664952642947
$ws.Range("A23").Value = "ARMA_I(0,7,0)"
$ws.Range("H23").Value = 97.25361541597267
$ws.Range("A24").Value = "ARMA_I(0,7,1)"
$ws.Range("H24").Value = 97.33924396059854
$ws.Range("A25").Value = "ARMA_I(0,7,2)"
$ws.Range("H25").Value = 97.28763606178265
$ws.Range("A26").Value = "ARMA_I(1,1,0)"
$ws.Range("B26").Value = 88.85368895367749
$ws.Range("A27").Value = "ARMA_I(1,1,1)"
$ws.Range("B27").Value = 91.8468054008783
$ws.Range("A28").Value = "ARMA_I(1,10,0)"
$ws.Range("I28").Value = 95.91481121943517
$ws.Range("A29").Value = "ARMA_I(1,10,1)"
$ws.Range("I29").Value = 95.87691392297177
$ws.Range("A30").Value = "ARMA_I(1,2,0)"
$ws.Range("C30").Value = 97.7092903514183
$ws.Range("A31").Value = "ARMA_I(1,2,1)"
$ws.Range("C31").Value = 97.88577715113331
$ws.Range("A32").Value = "ARMA_I(1,3,0)"
$ws.Range("D32").Value = 99.28963624262335
$ws.Range("A33").Value = "ARMA_I(1,3,1)"
$ws.Range("D33").Value = 99.26375222298115
$ws.Range("A34").Value = "ARMA_I(1,4,0)"
$ws.Range("E34").Value = 98.78417761786724
$ws.Range("A35").Value = "ARMA_I(1,4,1)"
$ws.Range("E35").Value = 98.74528000914586
$ws.Range("A36").Value = "ARMA_I(1,5,0)"
$ws.Range("F36").Value = 98.25490596901302
$ws.Range("A37").Value = "ARMA_I(1,5,1)"
$ws.Range("F37").Value = 98.32839331878446
$ws.Range("A38").Value = "ARMA_I(1,6,0)"
$ws.Range("G38").Value = 97.7971539362587
$ws.Range("A39").Value = "ARMA_I(1,6,1)"
$ws.Range("G39").Value = 97.73701677958006
$ws.Range("A40").Value = "ARMA_I(1,7,0)"
$ws.Range("H40").Value = 97.22558774092572
$ws.Range("A41").Value = "ARMA_I(1,7,1)"
$ws.Range("H41").Value = 97.24699659693512
$ws.Range("A42").Value = "ARMA_I(2,1,0)"
$ws.Range("B42").Value = 90.53224334979524
$ws.Range("A43").Value = "ARMA_I(2,1,2)"
$ws.Range("B43").Value = 92.61237230484846
$ws.Range("A44").Value = "ARMA_I(2,10,0)"
$ws.Range("I44").Value = 95.8771402532327
$ws.Range("A45").Value = "ARMA_I(2,10,2)"
$ws.Range("I45").Value = 95.81760157069749
$ws.Range("A46").Value = "ARMA_I(2,2,0)"
$ws.Range("C46").Value = 98.46355862844388
$ws.Range("A47").Value = "ARMA_I(2,2,2)"
$ws.Range("C47").Value = 98.1154257567485
$ws.Range("A48").Value = "ARMA_I(2,3,0)"
$ws.Range("D48").Value = 99.2903186757181
$ws.Range("A49").Value = "ARMA_I(2,3,2)"
$ws.Range("D49").Value = 99.23788476988348
$ws.Range("A50").Value = "ARMA_I(2,4,0)"
$ws.Range("E50").Value = 98.73713926589632
$ws.Range("A51").Value = "ARMA_I(2,4,2)"
$ws.Range("E51").Value = 98.67393416462255
$ws.Range("A52").Value = "ARMA_I(2,5,0)"
$ws.Range("F52").Value = 98.29375315960617
$ws.Range("A53").Value = "ARMA_I(2,5,2)"
$ws.Range("F53").Value = 98.16132394066372
$ws.Range("A54").Value = "ARMA_I(2,6,0)"
$ws.Range("G54").Value = 97.77557490282635
$ws.Range("A55").Value = "ARMA_I(2,6,2)"
$ws.Range("G55").Value = 97.70591282608083
$ws.Range("A56").Value = "ARMA_I(2,7,0)"
$ws.Range("H56").Value = 97.23947772932038
$ws.Range("A57").Value = "ARMA_I(2,7,2)"
$ws.Range("H57").Value = 97.34820116697385

# Apply header/label style (bold, centered, thin border) to newly introduced cells
$ws.Range("B1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("A51:A57").PasteSpecial(-4122)
$excel.CutCopyMode = 0
